{"js": "// Rename the merge-field placeholders used by this correspondence template\n// to the camelCase field names used by the Core DEV merge configuration:\n//   ${Date}                          -> ${currentDate}\n//   ${Employee Name}  (both usages)  -> ${personFirstName}\n//   ${Case Number}: ${Case Title}    -> ${caseNumber}: ${caseTitle}\nconst body = context.document.body;\n\n// ${Date} -> ${currentDate}\nlet dateResults = body.search(\"${Date}\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (const item of dateResults.items) {\n  item.insertText(\"${currentDate}\", \"Replace\");\n}\nawait context.sync();\n\n// ${Employee Name} -> ${personFirstName} (appears twice: salutation + body)\nlet nameResults = body.search(\"${Employee Name}\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\nfor (const item of nameResults.items) {\n  item.insertText(\"${personFirstName}\", \"Replace\");\n}\nawait context.sync();\n\n// ${Case Number}: ${Case Title} -> ${caseNumber}: ${caseTitle}\nlet caseResults = body.search(\"${Case Number}: ${Case Title}\", { matchCase: true });\ncaseResults.load(\"items\");\nawait context.sync();\nfor (const item of caseResults.items) {\n  item.insertText(\"${caseNumber}: ${caseTitle}\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Rename the merge-field placeholders used by this correspondence template\n# to the camelCase field names used by the Core DEV merge configuration:\n#   ${Date}                          -> ${currentDate}\n#   ${Employee Name}  (both usages)  -> ${personFirstName}\n#   ${Case Number}: ${Case Title}    -> ${caseNumber}: ${caseTitle}\n$d = $word.ActiveDocument\n\n# ${Date} -> ${currentDate}\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n[void]$find1.Find.Execute(\"`${Date}\", $true, $false, $false, $false, $false, $true, 1, $false, \"`${currentDate}\", 2)\n\n# ${Employee Name} -> ${personFirstName} (salutation + body paragraph)\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n[void]$find2.Find.Execute(\"Employee Name\", $true, $false, $false, $false, $false, $true, 1, $false, \"personFirstName\", 2)\n\n# ${Case Number}: ${Case Title} -> ${caseNumber}: ${caseTitle}\n$find3 = $d.Content\n$find3.Find.ClearFormatting()\n[void]$find3.Find.Execute(\"Case Number}: `${Case Title\", $true, $false, $false, $false, $false, $true, 1, $false, \"caseNumber}: `${caseTitle\", 2)\n"}
